$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "group1"
$ws.Range("B1").Value = "group2"
$ws.Range("C1").Value = "meandiff"
$ws.Range("D1").Value = "p-adj"
$ws.Range("E1").Value = "lower"
$ws.Range("F1").Value = "upper"
$ws.Range("G1").Value = "reject"

# Copy header style (bold font, thin borders, center/top alignment) from A1 to B1:G1
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2)
$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"
$ws.Range("C2").Value = -1.2765
$ws.Range("D2").Value = 0.001
$ws.Range("E2").Value = -1.7855
$ws.Range("F2").Value = -0.7675999999999999
$ws.Range("G2").Value = $true
